# Apply the documented change by editing the document's underlying OOXML.
# We read the package-wide WordOpenXML, perform precise textual substitutions
# that mirror the unified diff exactly, then push the corrected XML back into
# the document via Range.InsertXML (which replaces the content it targets).

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

# --- Change 1: insert a new (blank-ish) paragraph right after the title
#     paragraph and before the "احراز هویت..." Heading1 paragraph. The new
#     paragraph contains a single space run in the "B Nazanin" CS font at
#     size 24 (half-points: sz/szCs = 24 => 12pt).
$old1 = @'
<w:t>سایت نوبت دهی</w:t></w:r></w:p><w:p w14:paraId="53B37D58"
'@
$new1 = @'
<w:t>سایت نوبت دهی</w:t></w:r></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:cs="B Nazanin"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p w14:paraId="53B37D58"
'@
if ($xml.IndexOf($old1) -lt 0) {
    throw "anchor1 not found"
}
$xml = $xml.Replace($old1, $new1)

# --- Change 2: on the "پزشک کیست و چه ویژگی هایی دارد؟" Heading1 paragraph,
#     drop the stray <w:rFonts w:hint="cs"/> from the paragraph-mark's rPr
#     (the run's own rPr keeps it), then insert four new body paragraphs
#     immediately after it (before "نوبت دهی چیست و ویژگی های آن").
$old2 = @'
<w:p w14:paraId="3FD69786" w14:textId="646DCFD4" w:rsidR="00582D2C" w:rsidRDefault="00582D2C" w:rsidP="00582D2C"><w:pPr><w:pStyle w:val="Heading1"/><w:bidi/><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>پزشک کیست و چه ویژگی هایی دارد؟</w:t></w:r></w:p>
'@
$new2 = @'
<w:p w14:paraId="3FD69786" w14:textId="646DCFD4" w:rsidR="00582D2C" w:rsidRDefault="00582D2C" w:rsidP="00582D2C"><w:pPr><w:pStyle w:val="Heading1"/><w:bidi/><w:rPr><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>پزشک کیست و چه ویژگی هایی دارد؟</w:t></w:r></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>پزشک باید بتواند</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>زمان های خود را انتخاب کند.</w:t></w:r></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>پزشکی را درنظر بگیرید که از ساعت 8:00 تا 13:00 شیفت صبح تا ظهر و از ساعت 14:00 تا 21:00 شیفت بعد از ظهر تا شب کار میکند که به بیماران خود نوبت 30 دقیقه ای میدهد.</w:t></w:r></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>در صورتی که پزشک مثلا برای دوهفته به علت سفر نتواند نوبت ده</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی انجام دهد،</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> باید به کاربر اعلام شود.</w:t></w:r></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>پزشک باید لیستی از بیمارانی که با آن نوبت دارند داشته باشد و این لیست شامل(نام و نام خانوادگی بیمار، کدملی، ساعت شروع و پایان نوبت، تاریخ نوبت) باشد.</w:t></w:r></w:p>
'@
if ($xml.IndexOf($old2) -lt 0) {
    throw "anchor2 not found"
}
$xml = $xml.Replace($old2, $new2)

$result = $d.Content.InsertXML($xml)
Write-Output "applied"
